$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-5 (new sensor readings)
$rowData = @(
    @(45174.50694444445,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(45174.51388888889,18.764,14.012,0.672,41.453,35.05,14.973,48.024,22.737,10.338,15.906,18.091,17.396,5.623,14.883,21.254,12.147,0.253,0.891,220.953,41.013,13.718,27.887,14.504,2.901,24.806,12.186,10.767,12.876,18.499,1.533,43.923,8.355,16.923),
    @(45174.52083333334,15.103,11.294,0.472,33.301,27.908,12.023,46.976,18.302,8.343999999999999,12.614,14.118,13.988,4.363,11.89,17.14,9.819000000000001,0.169,0.67,175.494,33.195,10.966,22.578,11.707,2.123,22.956,9.782999999999999,8.662000000000001,10.31,14.636,0.861,43.315,6.597,13.623),
    @(45174.52777777778,21.91,16.42,0.6899999999999999,48,39.93,17.35,66.31999999999999,26.54,12.02,18.02,19.76,20.26,5.92,17.17,24.7,14.24,0.21,0.82,255.71,47.99,15.85,32.59,17.01,2.65,32.46,14.09,12.43,14.7,20.66,0.6,60.53,9.289999999999999,19.77)
)

$startRow = 2
for ($r = 0; $r -lt $rowData.Length; $r++) {
    $rowVals = $rowData[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowVals[$c]
    }
}

# Remove row 6 entirely (data series now only has 4 rows after the header)
$ws.Rows.Item(6).Delete()

# Widen columns C (3), J (10), and AB (28) from 7 to 8 characters
$ws.Range("C1").EntireColumn.ColumnWidth = 7.17
$ws.Range("J1").EntireColumn.ColumnWidth = 7.17
$ws.Range("AB1").EntireColumn.ColumnWidth = 7.17
